$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Stage format stamps for each style index (1-8) in a scratch area (column ZZ) ----
$ws.Range("A1").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)
$ws.Range("A35").Copy()
$ws.Range("ZZ2").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("ZZ3").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("ZZ4").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("ZZ5").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("ZZ6").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("ZZ7").PasteSpecial(-4122)
$ws.Range("ZZ7").WrapText = $true
$ws.Range("A35").Copy()
$ws.Range("ZZ8").PasteSpecial(-4122)
$ws.Range("ZZ8").WrapText = $true
$ws.Application.CutCopyMode = $false

# ---- Apply per-row format stamps + set row heights ----
# row 1
$ws.Range("ZZ1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Rows("1:1").RowHeight = 42
$ws.Range("A1").Value = 'SRS ID'
$ws.Range("B1").Value = 'Test case ID'

# row 2
$ws.Range("ZZ4").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Rows("2:2").RowHeight = 23.25
$ws.Range("A2").Value = 'SRS_Register_001'
$ws.Range("B2").Value = 'TC_Register_001'

# row 3
$ws.Range("ZZ4").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Rows("3:3").RowHeight = 23.25
$ws.Range("A3").Value = 'SRS_Register_002'
$ws.Range("B3").Value = 'TC_Register_002'

# row 4
$ws.Range("ZZ4").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Rows("4:4").RowHeight = 23.25
$ws.Range("A4").Value = 'SRS_Register_003'
$ws.Range("B4").Value = 'TC_Register_003'

# row 5
$ws.Range("ZZ4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Rows("5:5").RowHeight = 23.25
$ws.Range("A5").Value = 'SRS_Register_004'
$ws.Range("B5").Value = 'TC_Register_004'

# row 6
$ws.Range("ZZ4").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Rows("6:6").RowHeight = 23.25
$ws.Range("A6").Value = 'SRS_Register_005'
$ws.Range("B6").Value = 'TC_Register_005'

# row 7
$ws.Range("ZZ4").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Rows("7:7").RowHeight = 23.25
$ws.Range("A7").Value = 'SRS_Register_006'
$ws.Range("B7").Value = 'TC_Register_006'

# row 8
$ws.Range("ZZ4").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Rows("8:8").RowHeight = 23.25
$ws.Range("A8").Value = 'SRS_Register_007'
$ws.Range("B8").Value = 'TC_Register_007'

# row 9
$ws.Range("ZZ4").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Rows("9:9").RowHeight = 23.25
$ws.Range("A9").Value = 'SRS_Register_008'
$ws.Range("B9").Value = 'TC_Register_008'

# row 10
$ws.Range("ZZ4").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Rows("10:10").RowHeight = 23.25
$ws.Range("A10").Value = 'SRS_Register_009'
$ws.Range("B10").Value = 'TC_Register_009'

# row 11
$ws.Range("ZZ4").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Rows("11:11").RowHeight = 23.25
$ws.Range("A11").Value = 'SRS_Register_010'
$ws.Range("B11").Value = 'TC_Register_010'

# row 12
$ws.Range("ZZ4").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Rows("12:12").RowHeight = 23.25
$ws.Range("A12").Value = 'SRS_Register_011'
$ws.Range("B12").Value = 'TC_Register_011'

# row 13
$ws.Range("ZZ4").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Rows("13:13").RowHeight = 23.25
$ws.Range("A13").Value = 'SRS_Register_012'
$ws.Range("B13").Value = 'TC_Register_012'

# row 14
$ws.Range("ZZ4").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Rows("14:14").RowHeight = 23.25
$ws.Range("A14").Value = 'SRS_Register_013'
$ws.Range("B14").Value = 'TC_Register_013'

# row 15
$ws.Range("ZZ4").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Rows("15:15").RowHeight = 23.25
$ws.Range("A15").Value = 'SRS_Register_014'
$ws.Range("B15").Value = 'TC_Register_014'

# row 16
$ws.Range("ZZ4").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Rows("16:16").RowHeight = 23.25
$ws.Range("A16").Value = 'SRS_Register_015'
$ws.Range("B16").Value = 'TC_Register_015'

# row 17
$ws.Range("ZZ4").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Rows("17:17").RowHeight = 23.25
$ws.Range("A17").Value = 'SRS_Register_016'
$ws.Range("B17").Value = 'TC_Register_016'

# row 18
$ws.Range("ZZ4").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Rows("18:18").RowHeight = 23.25
$ws.Range("A18").Value = 'SRS_Register_017'
$ws.Range("B18").Value = 'TC_Register_017'

# row 19
$ws.Range("ZZ4").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Rows("19:19").RowHeight = 23.25
$ws.Range("A19").Value = 'SRS_Register_018'
$ws.Range("B19").Value = 'TC_Register_001'

# row 20
$ws.Range("ZZ4").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Rows("20:20").RowHeight = 23.25
$ws.Range("A20").Value = 'SRS_Register_019'
$ws.Range("B20").Value = 'TC_Register_002'

# row 21
$ws.Range("ZZ5").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Rows("21:21").RowHeight = 23.25
$ws.Range("A21").Value = 'SRS_Login_001'
$ws.Range("B21").Value = 'TC_Login_001 '

# row 22
$ws.Range("ZZ5").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("ZZ6").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Rows("22:22").RowHeight = 23.25
$ws.Range("A22").Value = 'SRS_Login_002'
$ws.Range("B22").Value = 'TC_Login_001  '

# row 23
$ws.Range("ZZ5").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Rows("23:23").RowHeight = 23.25
$ws.Range("A23").Value = 'SRS_Login_003'
$ws.Range("B23").Value = 'TC_Login_014'

# row 24
$ws.Range("ZZ5").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("ZZ6").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Rows("24:24").RowHeight = 45.75
$ws.Range("A24").Value = 'SRS_Login_004'
$ws.Range("B24").Value = 'TC_Login_001  
TC-Login_008'

# row 25
$ws.Range("ZZ5").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("ZZ6").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Rows("25:25").RowHeight = 42
$ws.Range("A25").Value = 'SRS_Login_005'
$ws.Range("B25").Value = 'TC_Login_002
TC_Login_003
TC_Login_004
TC_Login_009
TC_Login_0010
TC_Login_0011
'

# row 26
$ws.Range("ZZ5").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Rows("26:26").RowHeight = 23.25
$ws.Range("A26").Value = 'SRS_Supplier_001'
$ws.Range("B26").Value = 'TC_Supplier_001'

# row 27
$ws.Range("ZZ5").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("ZZ7").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$ws.Rows("27:27").RowHeight = 23.25
$ws.Range("A27").Value = 'SRS_Supplier_002'
$ws.Range("B27").Value = 'TC_Supplier_002
TC_Supplier_003'

# row 28
$ws.Range("ZZ5").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$ws.Rows("28:28").RowHeight = 23.25
$ws.Range("A28").Value = 'SRS_Supplier_003'
$ws.Range("B28").Value = 'TC_Supplier_004'

# row 29
$ws.Range("ZZ5").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Rows("29:29").RowHeight = 23.25
$ws.Range("A29").Value = 'SRS_Supplier_004'
$ws.Range("B29").Value = 'TC_Supplier_005'

# row 30
$ws.Range("ZZ5").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Rows("30:30").RowHeight = 23.25
$ws.Range("A30").Value = 'SRS_Supplier_005'
$ws.Range("B30").Value = 'TC_Supplier_005'

# row 31
$ws.Range("ZZ5").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Rows("31:31").RowHeight = 23.25
$ws.Range("A31").Value = 'SRS_Supplier_006'
$ws.Range("B31").Value = 'TC_Supplier_006'

# row 32
$ws.Range("ZZ5").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Rows("32:32").RowHeight = 23.25
$ws.Range("A32").Value = 'SRS_Supplier_007'
$ws.Range("B32").Value = 'TC_Supplier_007'

# row 33
$ws.Range("ZZ5").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Rows("33:33").RowHeight = 23.25
$ws.Range("A33").Value = 'SRS_Supplier_008'
$ws.Range("B33").Value = 'TC_Supplier_008'

# row 34
$ws.Range("ZZ5").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Rows("34:34").RowHeight = 23.25
$ws.Range("A34").Value = 'SRS_Supplier_009'
$ws.Range("B34").Value = 'TC_Supplier_009'

# row 35
$ws.Range("ZZ5").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("ZZ8").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Rows("35:35").RowHeight = 76.5
$ws.Range("A35").Value = 'SRS_Supplier_010'
$ws.Range("B35").Value = 'TC_Supplier_003
TC_Supplier_005
TC_Supplier_006
TC_Supplier_007'

# rows 36-38: blank, A=style5 B=style3 h=23.25 custom=True
$ws.Range("ZZ5").Copy()
$ws.Range("A36:A38").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B36:B38").PasteSpecial(-4122)
$ws.Rows("36:38").RowHeight = 23.25

# rows 39-39: blank, A=style5 B=style3 h=18.75 custom=False
$ws.Range("ZZ5").Copy()
$ws.Range("A39:A39").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B39:B39").PasteSpecial(-4122)
$ws.Rows("39:39").RowHeight = 18.75

# rows 40-44: blank, A=style5 B=style5 h=18.75 custom=False
$ws.Range("ZZ5").Copy()
$ws.Range("A40:A44").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B40:B44").PasteSpecial(-4122)
$ws.Rows("40:44").RowHeight = 18.75

# rows 45-45: blank, A=style5 B=style3 h=18.75 custom=False
$ws.Range("ZZ5").Copy()
$ws.Range("A45:A45").PasteSpecial(-4122)
$ws.Range("ZZ3").Copy()
$ws.Range("B45:B45").PasteSpecial(-4122)
$ws.Rows("45:45").RowHeight = 18.75

# rows 46-49: blank, A=style5 B=style5 h=18.75 custom=False
$ws.Range("ZZ5").Copy()
$ws.Range("A46:A49").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B46:B49").PasteSpecial(-4122)
$ws.Rows("46:49").RowHeight = 18.75

# rows 50-53: blank, A=style5 B=style5 h=15.75 custom=True
$ws.Range("ZZ5").Copy()
$ws.Range("A50:A53").PasteSpecial(-4122)
$ws.Range("ZZ5").Copy()
$ws.Range("B50:B53").PasteSpecial(-4122)
$ws.Rows("50:53").RowHeight = 15.75

# rows 54-249: blank, A=style2 B=style2 h=15.75 custom=True
$ws.Range("ZZ2").Copy()
$ws.Range("A54:A249").PasteSpecial(-4122)
$ws.Range("ZZ2").Copy()
$ws.Range("B54:B249").PasteSpecial(-4122)
$ws.Rows("54:249").RowHeight = 15.75

$ws.Application.CutCopyMode = $false
# ---- Clean up scratch stamp area ----
$ws.Range("ZZ1:ZZ8").Clear()

$ws.Range("A1").Select()